$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maggie Burton's time spent (row 7, column B) is updated from "20h 30m" to "21h 45m"
$ws.Range("B7").Value = "21h 45m"
$ws.Range("B7").Select()
